# Week 3 DSL doc edit script
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: remove the _GoBack bookmark that sits after
# "be sub-divided into internal and external languages."
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# Change 2: append two new sentences to the end of the "Domain Scripting"
# section paragraph (after "...properly implemented.")
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("properly implemented.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(" That can reduce the costs to implement new features responsibilities are decoupled and specialists operate on each aspect of the problem.")

# ---------------------------------------------------------------------------
# Change 3: the trailing empty "Heading2" paragraph (right after the
# "Internal vs External" heading) becomes the start of several new body
# paragraphs of text.
# ---------------------------------------------------------------------------

# -- paragraph: tab + "An internal domain specific language ..." (Normal style)
$n = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($n)
$p.Style = "Normal"
$rng = $p.Range
$rng.Collapse(1)
$rng.InsertAfter([char]9 + "An internal domain specific language is embedded within the context of its parent general purpose language, while an external resides outside. ")

# -- new paragraph: "Internal languages often use creative tricks ..."
$rng = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Format.FirstLineIndent = 36
$rng = $p.Range
$rng.Collapse(1)
$rng.InsertAfter("Internal languages often use creative tricks to improve the readability of their language such as operator overloading, removing optional punctuation, and defining no/op bubble words. The proposed ancestry query language could have implemented in C++ by overloading the GreaterThan and IndexInto operators. ")

# -- new paragraph: "There are limits to this approach ..."
$rng = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Format.FirstLineIndent = 36
$rng = $p.Range
$rng.Collapse(1)
$rng.InsertAfter("There are limits to this approach and challenges to internal languages. The business unit notices that requesting great grandparents appears in a large set of queries. They want to allow specifying the hierarchical levels as a sequence of equal signs, such that ")

$pStart = $p.Range.Start
$codeStart = $rng.End
$rng.InsertAfter("(me) > (parent) > (grandparent) > (great grandparent) > (cousins)")
$codeEnd = $pStart + ($codeStart - $pStart) + ("(me) > (parent) > (grandparent) > (great grandparent) > (cousins)").Length
$codeRng = $d.Range($codeStart, $codeEnd)
$codeRng.Font.Name = "Courier New"
$codeRng.Font.Size = 10

$rng = $p.Range
$rng.Collapse(0)
$rng.InsertAfter(" is equal to ")

$rng = $p.Range
$rng.Collapse(0)
$code2Start = $rng.Start
$rng.InsertAfter("(me) ===> (cousin)")
$code2End = $code2Start + ("(me) ===> (cousin)").Length
$code2Rng = $d.Range($code2Start, $code2End)
$code2Rng.Font.Name = "Courier New"
$code2Rng.Font.Size = 10

$rng = $p.Range
$rng.Collapse(0)
$rng.InsertAfter(". If the internal language does not expose an ===> then it cannot be overridden, and the feature cannot be implemented.")

# -- new paragraph: "To gain additional flexibility ..."
$rng = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Format.FirstLineIndent = 36
$rng = $p.Range
$rng.Collapse(1)
$rng.InsertAfter("To gain additional flexibility the development team needs to use define an external DSL language and parse the commands into an abstract representation. This flexibility comes at the cost of being more effort to maintain custom grammar files. ")

# -- new paragraph: "Common criticism of external DSL ..." + _GoBack bookmark
$rng = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Format.FirstLineIndent = 36
$rng = $p.Range
$rng.Collapse(1)
$rng.InsertAfter("Common criticism of external DSL is having to learn dozens of micro languages adds to the complexity of the system. This argument neglects to account of the sunk cost of leveraging the base API. The users need to understand this in some sh")

$rng = $p.Range
$rng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rng) | Out-Null

Write-Output "done-stage-2"
